# Auto-generated edit script: updates Leve profit/price calculation
# columns (H-N) on several worksheets to reflect refreshed market data.
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H40").Value = 12074.1
$wsALC.Range("I40").Value = 16480.143
$wsALC.Range("K40").Value = 16480.143
$wsALC.Range("M40").Value = -16305.143
$wsALC.Range("H70").Value = 1566
$wsALC.Range("I70").Value = 1431.25
$wsALC.Range("J70").Value = 1629.4117
$wsALC.Range("K70").Value = 4293.75
$wsALC.Range("L70").Value = 4888.2351
$wsALC.Range("M70").Value = -4023.75
$wsALC.Range("N70").Value = -5428.2351
$wsALC.Range("H73").Value = 1566
$wsALC.Range("I73").Value = 1431.25
$wsALC.Range("J73").Value = 1629.4117
$wsALC.Range("K73").Value = 4293.75
$wsALC.Range("L73").Value = 4888.2351
$wsALC.Range("M73").Value = -3357.75
$wsALC.Range("N73").Value = -6760.2351
$wsALC.Range("H137").Value = 2490434.2
$wsALC.Range("I137").Value = 4816999
$wsALC.Range("J137").Value = 8765.666999999999
$wsALC.Range("K137").Value = 14450997
$wsALC.Range("L137").Value = 26297.001
$wsALC.Range("M137").Value = -14448447
$wsALC.Range("N137").Value = -31397.001

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H31").Value = 9466.083000000001
$wsARM.Range("I31").Value = 2127.5715
$wsARM.Range("J31").Value = 19740
$wsARM.Range("K31").Value = 2127.5715
$wsARM.Range("L31").Value = 19740
$wsARM.Range("M31").Value = -1833.5715
$wsARM.Range("N31").Value = -20328
$wsARM.Range("H32").Value = 9174.152
$wsARM.Range("I32").Value = 8147.178
$wsARM.Range("K32").Value = 8147.178
$wsARM.Range("M32").Value = -7860.178
$wsARM.Range("H45").Value = 1815.7567
$wsARM.Range("I45").Value = 1650.6897
$wsARM.Range("J45").Value = 2414.125
$wsARM.Range("K45").Value = 1650.6897
$wsARM.Range("L45").Value = 2414.125
$wsARM.Range("M45").Value = -1273.6897
$wsARM.Range("N45").Value = -3168.125
$wsARM.Range("H74").Value = 1495.8591
$wsARM.Range("I74").Value = 1242.3889
$wsARM.Range("J74").Value = 2301
$wsARM.Range("K74").Value = 1242.3889
$wsARM.Range("L74").Value = 2301
$wsARM.Range("M74").Value = -368.3888999999999
$wsARM.Range("N74").Value = -4049
$wsARM.Range("H77").Value = 1495.8591
$wsARM.Range("I77").Value = 1242.3889
$wsARM.Range("J77").Value = 2301
$wsARM.Range("K77").Value = 6211.9445
$wsARM.Range("L77").Value = 11505
$wsARM.Range("M77").Value = -1843.9445
$wsARM.Range("N77").Value = -20241
$wsARM.Range("H122").Value = 1795.6471
$wsARM.Range("I122").Value = 1567.6666
$wsARM.Range("J122").Value = 2342.8
$wsARM.Range("K122").Value = 4702.9998
$wsARM.Range("L122").Value = 7028.400000000001
$wsARM.Range("M122").Value = -2252.9998
$wsARM.Range("N122").Value = -11928.4
$wsARM.Range("H132").Value = 11906402
$wsARM.Range("I132").Value = 14707114
$wsARM.Range("K132").Value = 44121342
$wsARM.Range("M132").Value = -44118812

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H21").Value = 19987.5
$wsBSM.Range("J21").Value = 19987.5
$wsBSM.Range("L21").Value = 19987.5
$wsBSM.Range("N21").Value = -20459.5
$wsBSM.Range("H26").Value = 12920.615
$wsBSM.Range("I26").Value = 6852.5713
$wsBSM.Range("K26").Value = 6852.5713
$wsBSM.Range("M26").Value = -6560.5713
$wsBSM.Range("H28").Value = 19999.857
$wsBSM.Range("I28").Value = 19999.5
$wsBSM.Range("K28").Value = 19999.5
$wsBSM.Range("M28").Value = -19705.5
$wsBSM.Range("H134").Value = 2020.122
$wsBSM.Range("I134").Value = 1486.1177
$wsBSM.Range("K134").Value = 4458.3531
$wsBSM.Range("M134").Value = -1923.3531

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H11").Value = 1550.122
$wsCUL.Range("I11").Value = 1606.5385
$wsCUL.Range("J11").Value = 450
$wsCUL.Range("K11").Value = 4819.6155
$wsCUL.Range("L11").Value = 1350
$wsCUL.Range("M11").Value = -4679.6155
$wsCUL.Range("N11").Value = -1630
$wsCUL.Range("H22").Value = 17015.46
$wsCUL.Range("I22").Value = 22270.2
$wsCUL.Range("J22").Value = 13731.25
$wsCUL.Range("K22").Value = 66810.60000000001
$wsCUL.Range("L22").Value = 41193.75
$wsCUL.Range("M22").Value = -66641.60000000001
$wsCUL.Range("N22").Value = -41531.75
$wsCUL.Range("H27").Value = 17015.46
$wsCUL.Range("I27").Value = 22270.2
$wsCUL.Range("J27").Value = 13731.25
$wsCUL.Range("K27").Value = 66810.60000000001
$wsCUL.Range("L27").Value = 41193.75
$wsCUL.Range("M27").Value = -66708.60000000001
$wsCUL.Range("N27").Value = -41397.75
$wsCUL.Range("H107").Value = 5388.2954
$wsCUL.Range("J107").Value = 7116.778
$wsCUL.Range("L107").Value = 21350.334
$wsCUL.Range("N107").Value = -25190.334

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H16").Value = 3115.2222
$wsLTW.Range("I16").Value = 2962.2856
$wsLTW.Range("J16").Value = 3650.5
$wsLTW.Range("K16").Value = 2962.2856
$wsLTW.Range("L16").Value = 3650.5
$wsLTW.Range("M16").Value = -2792.2856
$wsLTW.Range("N16").Value = -3990.5
$wsLTW.Range("H46").Value = 3054.85
$wsLTW.Range("I46").Value = 831.0476
$wsLTW.Range("J46").Value = 5512.737
$wsLTW.Range("K46").Value = 831.0476
$wsLTW.Range("L46").Value = 5512.737
$wsLTW.Range("M46").Value = -643.0476
$wsLTW.Range("N46").Value = -5888.737
$wsLTW.Range("H68").Value = 3574.75
$wsLTW.Range("I68").Value = 3476.923
$wsLTW.Range("J68").Value = 3659.5334
$wsLTW.Range("K68").Value = 3476.923
$wsLTW.Range("L68").Value = 3659.5334
$wsLTW.Range("M68").Value = -2727.923
$wsLTW.Range("N68").Value = -5157.5334
$wsLTW.Range("H71").Value = 3574.75
$wsLTW.Range("I71").Value = 3476.923
$wsLTW.Range("J71").Value = 3659.5334
$wsLTW.Range("K71").Value = 17384.615
$wsLTW.Range("L71").Value = 18297.667
$wsLTW.Range("M71").Value = -13640.615
$wsLTW.Range("N71").Value = -25785.667
$wsLTW.Range("H132").Value = 3946.0588
$wsLTW.Range("I132").Value = 3504.4119
$wsLTW.Range("J132").Value = 4387.706
$wsLTW.Range("K132").Value = 10513.2357
$wsLTW.Range("L132").Value = 13163.118
$wsLTW.Range("M132").Value = -7983.235700000001
$wsLTW.Range("N132").Value = -18223.118
$wsLTW.Range("H141").Value = 28811.334
$wsLTW.Range("J141").Value = 28811.334
$wsLTW.Range("L141").Value = 28811.334
$wsLTW.Range("N141").Value = -39171.334

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H81").Value = 4540
$wsWVR.Range("I81").Value = 1350
$wsWVR.Range("J81").Value = 6666.6665
$wsWVR.Range("K81").Value = 2700
$wsWVR.Range("L81").Value = 13333.333
$wsWVR.Range("M81").Value = -1639
$wsWVR.Range("N81").Value = -15455.333
$wsWVR.Range("H84").Value = 4540
$wsWVR.Range("I84").Value = 1350
$wsWVR.Range("J84").Value = 6666.6665
$wsWVR.Range("K84").Value = 13500
$wsWVR.Range("L84").Value = 66666.66500000001
$wsWVR.Range("M84").Value = -8196
$wsWVR.Range("N84").Value = -77274.66500000001
$wsWVR.Range("H136").Value = 467420.34
$wsWVR.Range("I136").Value = 530914.2
$wsWVR.Range("J136").Value = 1798.8334
$wsWVR.Range("K136").Value = 1592742.6
$wsWVR.Range("L136").Value = 5396.5002
$wsWVR.Range("M136").Value = -1590192.6
$wsWVR.Range("N136").Value = -10496.5002

Write-Host "Updated 169 cells across 6 worksheets (ALC, ARM, BSM, CUL, LTW, WVR)"
